$d = $word.ActiveDocument

function Insert-ItemAfter($anchorText, $wholeWord, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $wholeWord, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor text not found: $anchorText"
    }
    $para = $rng.Paragraphs(1)
    $para.Range.InsertParagraphAfter()
    $newPara = $para.Next()
    $newPara.Range.Text = $newText
    $newPara.Range.ListFormat.ListLevelNumber = 2
    return $newPara
}

# 1) After "Menu sends active reset signal to other modules when appropriate"
Insert-ItemAfter "Menu sends active reset signal to other modules when appropriate" $false "Optional: Maintain a high score table for each song" | Out-Null

# 2) After "Musical Score Loader" add four sub-items (in order)
$anchor = "Musical Score Loader"
$items = @(
    "Song files are properly stored / accessible from EEPROM",
    "A single song is loadable and does not have any invalid output",
    "All available songs load and play correctly to their own tempos",
    "Extraordinarily optional: be able to feed in a MIDI file to play"
)
foreach ($item in $items) {
    $newPara = Insert-ItemAfter $anchor $true $item
    $anchor = $item
}

# 3) After "Display" (standalone heading) add sub-items, with the testbench one using InsertXML
$rng = $d.Content
$found = $rng.Find.Execute("Display", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Display anchor not found" }
$para = $rng.Paragraphs(1)
$para.Range.InsertParagraphAfter()
$newPara = $para.Next()
$newPara.Range.Text = "A single note blob moves across the screen properly (ease in right, ease out left)"
$newPara.Range.ListFormat.ListLevelNumber = 2

$anchor = "A single note blob moves across the screen properly (ease in right, ease out left)"
$newPara = Insert-ItemAfter $anchor $false "The cstringdisp module is integrated and shows the score, current pitch"

# testbench paragraph with proofErr spell-check markers around "testbench"
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("The cstringdisp module is integrated and shows the score, current pitch", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "cstringdisp anchor not found" }
$para2 = $rng2.Paragraphs(1)
$para2.Range.InsertParagraphAfter()
$tbPara = $para2.Next()
$tbPara.Range.ListFormat.ListLevelNumber = 2
$tbXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t xml:space='preserve'>Creating a </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>testbench</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> to simulate inputs from the game logic, test hit pitches</w:t></w:r></w:p>"
$tbPara.Range.InsertXML($tbXml)

$anchor = "Creating a testbench to simulate inputs from the game logic, test hit pitches"
$items2 = @(
    "All note blobs are onscreen and transition smoothly",
    "Optional: Load background images in, such as a recorder finger chart",
    "Optional: Use bitmaps instead of notes",
    "Optional: Cool effects like fading notes and changing colors"
)
foreach ($item in $items2) {
    $newPara = Insert-ItemAfter $anchor $false $item
    $anchor = $item
}

# 4) Fix "Integration with FFT complete, can play an entire song" -- merge the two runs into one
$found3 = $d.Content.Find.Execute("Integration with FFT complete, can play an entire song", $true, $false, $false, $false, $false, $true, 1, $false, "Integration with FFT complete, can play an entire song", 2)
if (-not $found3) {
    throw "Integration with FFT text not found"
}

# 5) Add an extra empty paragraph at the very end (before sectPr)
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertParagraphAfter()

# 6) Move the _GoBack bookmark into the middle of "scale" in the "Whole system..." paragraph
$bookmarks = $d.Bookmarks
if ($bookmarks.Exists("_GoBack")) {
    $bookmarks.Item("_GoBack").Delete()
}
$rng4 = $d.Content
$found4 = $rng4.Find.Execute("Whole system responds well to at least one instrument for all C major s", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found4) { throw "Whole system anchor not found" }
$bmPos = $rng4.End
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output "DONE"
